$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3935.6875
$ws.Range("I40").Value = 1424.25
$ws.Range("J40").Value = 4294.4644
$ws.Range("K40").Value = 1424.25
$ws.Range("L40").Value = 4294.4644
$ws.Range("M40").Value = -1249.25
$ws.Range("N40").Value = -4644.4644
# Row 69
$ws.Range("H69").Value = 47649892
$ws.Range("I69").Value = 253006.5
$ws.Range("K69").Value = 759019.5
$ws.Range("M69").Value = -758145.5
# Row 72
$ws.Range("H72").Value = 47649892
$ws.Range("I72").Value = 253006.5
$ws.Range("K72").Value = 2277058.5
$ws.Range("M72").Value = -2272690.5
# Row 76
$ws.Range("H76").Value = 2062037.6
$ws.Range("I76").Value = 3272114.2
$ws.Range("J76").Value = 4907.4
$ws.Range("K76").Value = 3272114.2
$ws.Range("L76").Value = 4907.4
$ws.Range("M76").Value = -3271799.2
$ws.Range("N76").Value = -5537.4
# Row 79
$ws.Range("H79").Value = 2062037.6
$ws.Range("I79").Value = 3272114.2
$ws.Range("J79").Value = 4907.4
$ws.Range("K79").Value = 3272114.2
$ws.Range("L79").Value = 4907.4
$ws.Range("M79").Value = -3271022.2
$ws.Range("N79").Value = -7091.4
# Row 101
$ws.Range("H101").Value = 83335160
$ws.Range("I101").Value = 125000190
$ws.Range("J101").Value = 5092.5
$ws.Range("K101").Value = 375000570
$ws.Range("L101").Value = 15277.5
$ws.Range("M101").Value = -374998948
$ws.Range("N101").Value = -18521.5
# Row 112
$ws.Range("H112").Value = 5021.2144
$ws.Range("J112").Value = 5021.2144
$ws.Range("L112").Value = 15063.6432
$ws.Range("N112").Value = -17279.6432

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5486
$ws.Range("I32").Value = 3607.2195
$ws.Range("K32").Value = 3607.2195
$ws.Range("M32").Value = -3320.2195
# Row 34
$ws.Range("H34").Value = 11999
$ws.Range("I34").Value = 11999
$ws.Range("K34").Value = 11999
$ws.Range("M34").Value = -11728
# Row 63
$ws.Range("H63").Value = 3321.6667
$ws.Range("I63").Value = 3356.4285
$ws.Range("J63").Value = 3200
$ws.Range("K63").Value = 3356.4285
$ws.Range("L63").Value = 3200
$ws.Range("M63").Value = -2670.4285
$ws.Range("N63").Value = -4572
# Row 66
$ws.Range("H66").Value = 3321.6667
$ws.Range("I66").Value = 3356.4285
$ws.Range("J66").Value = 3200
$ws.Range("K66").Value = 16782.1425
$ws.Range("L66").Value = 16000
$ws.Range("M66").Value = -13350.1425
$ws.Range("N66").Value = -22864
# Row 74
$ws.Range("H74").Value = 214543.67
$ws.Range("J74").Value = 379255.44
$ws.Range("L74").Value = 379255.44
$ws.Range("N74").Value = -381003.44
# Row 77
$ws.Range("H77").Value = 214543.67
$ws.Range("J77").Value = 379255.44
$ws.Range("L77").Value = 1896277.2
$ws.Range("N77").Value = -1905013.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 6494959.5
$ws.Range("I107").Value = 7144255.5
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 7144255.5
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -7142335.5
$ws.Range("N107").Value = -5840

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 113338.766
$ws.Range("I132").Value = 73134
$ws.Range("K132").Value = 219402
$ws.Range("M132").Value = -216872

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 789.3333
$ws.Range("I3").Value = 789.3333
$ws.Range("K3").Value = 2367.9999
$ws.Range("M3").Value = -2255.9999
# Row 52
$ws.Range("H52").Value = 1050.6666
$ws.Range("J52").Value = 1050.6666
$ws.Range("L52").Value = 3151.9998
$ws.Range("N52").Value = -3683.9998
# Row 87
$ws.Range("H87").Value = 12212.857
$ws.Range("I87").Value = 8298
$ws.Range("K87").Value = 24894
$ws.Range("M87").Value = -23646
# Row 90
$ws.Range("H90").Value = 12212.857
$ws.Range("I90").Value = 8298
$ws.Range("K90").Value = 74682
$ws.Range("M90").Value = -68442

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 150
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -38
$ws.Range("N5").ClearContents()
# Row 11
$ws.Range("H11").Value = 10000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 10000
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -10278
# Row 14
$ws.Range("H14").Value = 997.7778
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 997.7778
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 997.7778
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1333.7778
# Row 63
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372
# Row 66
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864
# Row 102
$ws.Range("H102").Value = 9712251
$ws.Range("I102").Value = 22224422
$ws.Range("J102").Value = 2761045
$ws.Range("K102").Value = 22224422
$ws.Range("L102").Value = 2761045
$ws.Range("M102").Value = -22222800
$ws.Range("N102").Value = -2764289
# Row 113
$ws.Range("H113").Value = 55556892
$ws.Range("I113").Value = 83334090
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 83334090
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -83331920
$ws.Range("N113").Value = -6840
# Row 132
$ws.Range("H132").Value = 3767
$ws.Range("I132").Value = 3195.8572
$ws.Range("K132").Value = 9587.571599999999
$ws.Range("M132").Value = -7057.571599999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 89580.89999999999
$ws.Range("I22").Value = 178361.8
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 178361.8
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -178066.8
$ws.Range("N22").Value = -1390
# Row 27
$ws.Range("H27").Value = 89580.89999999999
$ws.Range("I27").Value = 178361.8
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 178361.8
$ws.Range("L27").Value = 800
$ws.Range("M27").Value = -178254.8
$ws.Range("N27").Value = -1014
# Row 46
$ws.Range("H46").Value = 5345.6895
$ws.Range("I46").Value = 4033.9412
$ws.Range("J46").Value = 7204
$ws.Range("K46").Value = 4033.9412
$ws.Range("L46").Value = 7204
$ws.Range("M46").Value = -3845.9412
$ws.Range("N46").Value = -7580

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# Row 31
$ws.Range("H31").Value = 16006.333
$ws.Range("J31").Value = 21509.5
$ws.Range("L31").Value = 21509.5
$ws.Range("N31").Value = -22205.5
# Row 54
$ws.Range("H54").Value = 20555.223
$ws.Range("J54").Value = 39998.5
$ws.Range("L54").Value = 39998.5
$ws.Range("N54").Value = -41038.5
# Row 62
$ws.Range("H62").Value = 5872.796
$ws.Range("J62").Value = 9125.772000000001
$ws.Range("L62").Value = 9125.772000000001
$ws.Range("N62").Value = -10373.772
# Row 64
$ws.Range("H64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40496
# Row 65
$ws.Range("H65").Value = 5872.796
$ws.Range("J65").Value = 9125.772000000001
$ws.Range("L65").Value = 45628.86
$ws.Range("N65").Value = -51868.86
# Row 67
$ws.Range("H67").Value = 40000
$ws.Range("J67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41716
# Row 81
$ws.Range("H81").Value = 18519528
$ws.Range("I81").Value = 27778542
$ws.Range("J81").Value = 1500.6666
$ws.Range("K81").Value = 55557084
$ws.Range("L81").Value = 3001.3332
$ws.Range("M81").Value = -55556023
$ws.Range("N81").Value = -5123.3332
# Row 84
$ws.Range("H84").Value = 18519528
$ws.Range("I84").Value = 27778542
$ws.Range("J84").Value = 1500.6666
$ws.Range("K84").Value = 277785420
$ws.Range("L84").Value = 15006.666
$ws.Range("M84").Value = -277780116
$ws.Range("N84").Value = -25614.666
# Row 107
$ws.Range("H107").Value = 62503970
$ws.Range("I107").Value = 71432870
$ws.Range("J107").Value = 1650.5
$ws.Range("K107").Value = 214298610
$ws.Range("L107").Value = 4951.5
$ws.Range("M107").Value = -214296690
$ws.Range("N107").Value = -8791.5
# Row 113
$ws.Range("H113").Value = 1047.8214
$ws.Range("I113").Value = 802.55554
$ws.Range("J113").Value = 1489.3
$ws.Range("K113").Value = 2407.66662
$ws.Range("L113").Value = 4467.9
$ws.Range("M113").Value = -237.66662
$ws.Range("N113").Value = -8807.9
